$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.064.24"
$ws.Range("E2").Value = "  -1.64%  "
$ws.Range("D3").Value = "2.105.38"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  -0.83%  "
$ws.Range("D5").Value = "'347.93"
$ws.Range("E5").Value = "  +3.21%  "
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("D7").Value = "'0.5169"
$ws.Range("E7").Value = "  -1.53%  "
$ws.Range("E8").Value = "  -2.64%  "
$ws.Range("D9").Value = "'52.25"
$ws.Range("E9").Value = "  -4.09%  "
$ws.Range("E10").Value = "  -1.99%  "
$ws.Range("D11").Value = "'1.172"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").Value = "'25.48"
$ws.Range("E12").Value = "  +3.74%  "
$ws.Range("D13").Value = "2.106.17"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").Value = "'8.245"
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").Value = "'6.732"
$ws.Range("E15").Value = "  -1.90%  "
$ws.Range("D16").Value = "'99.24"
$ws.Range("E16").Value = "  +2.21%  "
$ws.Range("E17").Value = "  -2.40%  "
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").Value = "'20.91"
$ws.Range("E19").Value = "  +7.66%  "
$ws.Range("D20").Value = "'0.06677"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("D22").Value = "'6.242"
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("D23").Value = "30.150.71"
$ws.Range("E23").Value = "  -1.59%  "
$ws.Range("D24").Value = "'12.72"
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("D25").Value = "'2.348"
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("D26").Value = "2.356.13"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").Value = "'21.99"
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("D28").Value = "'2.547"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "'162.52"
$ws.Range("E29").Value = "  -0.95%  "
$ws.Range("D30").Value = "'133.62"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").Value = "'1.176"
$ws.Range("E31").Value = "  -3.34%  "
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").Value = "'1.641"
$ws.Range("E33").Value = "  -0.41%  "
$ws.Range("D34").Value = "'6.238"
$ws.Range("D35").Value = "'3.961"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").Value = "'5.915"
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("E37").Value = "  -3.76%  "
$ws.Range("D38").Value = "'0.02575"
$ws.Range("E38").Value = "  -2.15%  "
$ws.Range("D39").Value = "'0.06797"
$ws.Range("E39").Value = "  -0.80%  "
$ws.Range("E40").Value = "  -0.99%  "
$ws.Range("D41").Value = "'12.56"
$ws.Range("D42").Value = "'0.6821"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("D43").Value = "'1.293"
$ws.Range("E43").Value = "  +2.64%  "
$ws.Range("D44").Value = "'14.29"
$ws.Range("E44").Value = "  -3.77%  "
$ws.Range("D45").Value = "'0.6385"
$ws.Range("E45").Value = "  -1.11%  "
$ws.Range("D46").Value = "'2.295"
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("D47").Value = "'0.00000000364"
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("E48").Value = "  -1.36%  "
$ws.Range("E49").Value = "  -2.76%  "
$ws.Range("D50").Value = "'82.52"
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("D51").Value = "'0.07231"
$ws.Range("E51").Value = "  +0.34%  "
